$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '57.765.05'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -6.41%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.898.10'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -4.91%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '552.10'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -5.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '121.51'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -6.75%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '2.892.05'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -5.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.494'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.125'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -9.50%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '4.74'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -9.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.433'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000211'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -9.04%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '31.48'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.99%  '
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.366.95'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.11%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.884.30'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -5.21%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.50'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '57.628.30'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.63%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '409.03'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -8.45%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.80'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.89%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.651'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.72'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -8.22%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.57'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '76.79'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.97%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  -3.94%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.91'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.39%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.12'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.11%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.03'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -6.80%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '24.59'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.00%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0947'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.02'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -13.24%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.898'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.70%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.34'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -6.28%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '48.36'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.36'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0₃0617'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -11.70%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0343'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -7.98%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.105'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.14%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '363.29'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.87%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.604.98'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.34'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -7.79%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  -4.27%  '
$ws.Range('E47').Value = '  -4.59%  '
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.93'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.81%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '22.20'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -6.71%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.95'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.25%  '
